$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, pushing existing rows 37:81 down to 38:82
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly record
$row = 37
$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44904
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100114007
$ws.Cells.Item($row, 7).Value = "Jengibre"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 400
$ws.Cells.Item($row, 11).Value = 14000
$ws.Cells.Item($row, 12).Value = 15000
$ws.Cells.Item($row, 13).Value = 14500
$ws.Cells.Item($row, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 1115
$ws.Cells.Item($row, 17).Value = 13
$ws.Cells.Item($row, 18).Value = "Hortaliza"
